$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) string values
$ws.Range("B1").Value = "metadata4ing_IRI"
$ws.Range("C1").Value = "metadata4ing_DESC"
$ws.Range("D1").Value = "MS_IRI"
$ws.Range("E1").Value = "MS_DESC"

# New header cell F1 - set value then copy format from E1 (same header style)
$ws.Range("F1").Value = "MS_DEF"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Update data row (row 2) string values
$ws.Range("B2").Value = "http://www.w3.org/ns/prov#Activity"
$ws.Range("C2").Value = "{'label': None, 'prefLabel': 'Activity', 'altLabel': None, 'name': 'Activity'}"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/NCIT_C43431"
$ws.Range("E2").Value = "{'label': 'Activity'}"

# New data cell F2 - set value then copy format from E2 (no special style)
$ws.Range("F2").Value = "[]"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
